$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-31 Monday" "2025-04-01 Tuesday"

Replace-Text "673÷2=" "531÷7="
Replace-Text "862÷8=" "772÷2="
Replace-Text "205÷3=" "719÷2="
Replace-Text "138÷3=" "443÷8="
Replace-Text "735÷4=" "847÷4="
Replace-Text "269÷4=" "105÷2="
Replace-Text "786÷5=" "847÷3="
Replace-Text "188÷2=" "555÷2="
Replace-Text "120÷9=" "614÷5="
Replace-Text "199÷7=" "855÷8="
Replace-Text "201÷2=" "834÷4="
Replace-Text "356÷9=" "200÷8="
Replace-Text "912÷7=" "959÷3="
Replace-Text "856÷4=" "361÷9="
Replace-Text "827÷3=" "520÷9="
Replace-Text "259÷3=" "714÷8="
Replace-Text "188÷4=" "695÷8="
Replace-Text "140÷9=" "144÷3="
Replace-Text "467÷8=" "144÷7="
Replace-Text "845÷4=" "469÷3="
Replace-Text "430÷4=" "273÷7="
Replace-Text "134÷7=" "573÷6="
Replace-Text "197÷6=" "166÷7="
Replace-Text "171÷3=" "351÷2="
Replace-Text "327÷2=" "886÷7="
